$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Random-from-ADC-on-start related BoM changes -----------------------
# R3 changes from a 68 Ohm resistor to a 1k Ohm resistor, and two new
# resistors (R4 = 68 Ohm, R5 = 1k Ohm) are introduced right after it.
# Q1/Q2 change from AO3404A n-channel MOSFETs to MMBT2222A NPN BJT
# transistors (used to seed a random value off analog noise at boot).

$bjtUrl  = "https://lcsc.com/product-detail/Bipolar-Transistors-BJT_PANJIT-International-MMBT2222A_C2992465.html"
$bjtPart = "MMBT2222A"
$bjtPkg  = "SOT-23"
$bjtDesc = "npn bjt transistor"

$resUrl68  = "https://lcsc.com/product-detail/Chip-Resistor-Surface-Mount_PANASONIC-ERJPA3J680V_C445804.html"
$resPart68 = "ERJPA3J680V"
$resDesc68 = "68 Ohm resistor"

$resUrl1k  = "https://lcsc.com/product-detail/Chip-Resistor-Surface-Mount_PANASONIC-ERJPA3J102V_C441891.html"
$resPart1k = "ERJPA3J102V"
$resDesc1k = "1k Ohm resistor"

# Insert two fresh rows after existing row 66 (the old R3 row) so the
# rows below (Q1, Q2, U1) shift down from 67/68/69 to 69/70/71.
$ws.Rows(67).Insert()
$ws.Rows(67).Insert()

# Row 69 (was 67): Q1 -> MMBT2222A transistor (was AO3404A mosfet)
$ws.Range("B69").Value = $bjtUrl
$ws.Range("E69").Value = $bjtDesc
$ws.Range("C69").Value = $bjtPart
$ws.Range("D69").Value = $bjtPkg

# Row 70 (was 68): Q2 -> MMBT2222A transistor (was AO3404A mosfet)
$ws.Range("B70").Value = $bjtUrl
$ws.Range("E70").Value = $bjtDesc
$ws.Range("C70").Value = $bjtPart
$ws.Range("D70").Value = $bjtPkg

# Row 67 (new): R4, 68 Ohm resistor
$ws.Range("A67").Value = "R4"
$ws.Range("B67").Value = $resUrl68
$ws.Range("C67").Value = $resPart68
$ws.Range("D67").Value = "0603"
$ws.Range("E67").Value = $resDesc68

# Row 68 (new): R5, 1k Ohm resistor
$ws.Range("A68").Value = "R5"
$ws.Range("B68").Value = $resUrl1k
$ws.Range("C68").Value = $resPart1k
$ws.Range("D68").Value = "0603"
$ws.Range("E68").Value = $resDesc1k

# Row 66: R3 -> now a 1k Ohm resistor (was 68 Ohm)
$ws.Range("B66").Value = $resUrl1k
$ws.Range("C66").Value = $resPart1k
$ws.Range("D66").Value = "0603"
$ws.Range("E66").Value = $resDesc1k

# Row 71 (was 69, U1/ATTINY10) is left untouched content-wise; it simply
# shifted down two rows from the inserts above.

# --- Hyperlink on B7 (the LED's LCSC product-detail URL) -----------------
$linkCell = $ws.Range("B7")
$ws.Hyperlinks.Add($linkCell, $linkCell.Value2)

# --- Remove the logo picture in the corner --------------------------------
if ($ws.Shapes.Count -gt 0) {
    $ws.Shapes.Item(1).Delete()
}

# --- Update selection / scroll position -----------------------------------
$ws.Range("B7").Select()
